# Add an "Address" column (new column F) in front of the existing
# "District" column, which shifts from F to G. The address text for each
# teacher is the school/location portion of column B (everything before the
# trailing district name), with the comma+space separators removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; existing column F ("District") shifts to G.
$ws.Columns("F:F").Insert()

# New header cells.
$ws.Range("F2").Value = "Address"
$ws.Range("F3").Value = ""

$addresses = @(
    @{Row=4; Addr='S S High School SomenahalliHosadurga'},
    @{Row=5; Addr='S R D H S RenukapuraChallakere'},
    @{Row=6; Addr=''},
    @{Row=7; Addr='S C K H S MudhureHosadurga'},
    @{Row=8; Addr='G H S ThippareddyhallyChallakere'},
    @{Row=9; Addr='G H S PilaliHiriyur'},
    @{Row=10; Addr='Govt. High School Kodagavallihatty'},
    @{Row=11; Addr='S M H School HanumanthanahallyChallakere'},
    @{Row=12; Addr='Maharani Composit P U College'},
    @{Row=13; Addr='G H S Vaddarasiddavvanahally'},
    @{Row=14; Addr='G J C Hiriyur'},
    @{Row=15; Addr='G H S Kurubarahally'},
    @{Row=16; Addr='S V S R H SChowlurChallakere'},
    @{Row=17; Addr='S V R High School KondapuraHosadurga'},
    @{Row=18; Addr='Govt. High School Ghatti HosahalliHolalkere'},
    @{Row=19; Addr='N S C B H S AdanurHolalkere'},
    @{Row=20; Addr='G H S GoolihattiHosadurga'},
    @{Row=21; Addr='G H S B G KereMolakalmuru'},
    @{Row=22; Addr='Bapuji High School Challakere'},
    @{Row=23; Addr='S M R H S BalenahalluChallakere'},
    @{Row=24; Addr='Sri Siddarameshwara'},
    @{Row=25; Addr='S P P R H SNannivalaChallakere'},
    @{Row=26; Addr='S J R High School AralahalliHosadurga'},
    @{Row=27; Addr='S K R H School NeralakereHosadurga'},
    @{Row=28; Addr='S K H S Bahaddurghatta'},
    @{Row=29; Addr='G H S NDevarahallyChallakere'},
    @{Row=30; Addr='G H S Laxmisagar'},
    @{Row=31; Addr='G H S MaskalHiriyur'},
    @{Row=32; Addr='G H S IkkanurHiriyur'},
    @{Row=33; Addr='Govt. High School Basavana Shivanakere'},
    @{Row=34; Addr='Govt. P U CollegeHigh School Section JavanagondanahallyHiriyur'},
    @{Row=35; Addr='G H S ObalapuraChallakere'},
    @{Row=36; Addr='G J C YaraballyHiriyur'},
    @{Row=37; Addr='G J C (HS) Molakalmuru'},
    @{Row=38; Addr='Govt. Composite P U College (High School Section)AvinahattiHolalkere'},
    @{Row=39; Addr='Govt. High School PagadalabandeChallakere'},
    @{Row=40; Addr='G H S VasanthanagaraHiriyur'},
    @{Row=41; Addr='S S M G H School Alagavadi'},
    @{Row=42; Addr='Sri M J R High School HoovinaholeHiriyur'},
    @{Row=43; Addr='A H S DummiHolalkere'},
    @{Row=44; Addr='G H S KaparahallyChallakere'},
    @{Row=45; Addr='R R Jr. CollegeMeerasabihallyChallakere'},
    @{Row=46; Addr='S P P H School KaluvehallyChallakere'},
    @{Row=47; Addr='Boys Jr. College'},
    @{Row=48; Addr='G H S G R Halli'},
    @{Row=49; Addr='Govt. High School ThammenahalliMolakalmuru'},
    @{Row=50; Addr='S P S R H School Molakalmuru'},
    @{Row=51; Addr='G J C Hosadurga'},
    @{Row=52; Addr=''}
)

foreach ($item in $addresses) {
    $ws.Cells.Item($item.Row, 6).Value = $item.Addr
}

